$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 141; this shifts rows 141:202 down to 142:203
$ws.Rows(141).Insert()

# Populate the new row 141 with the new record's data.
# Columns A, B, C, E, F, G, H, I, J, K are identical to the surrounding rows
# in this block (same market/region/category/variety), so copy them from row 142.
$ws.Range("A141").Value2 = $ws.Range("A142").Value2
$ws.Range("B141").Value2 = $ws.Range("B142").Value2
$ws.Range("C141").Value2 = $ws.Range("C142").Value2
$ws.Range("D141").Value2 = 45029
$ws.Range("E141").Value2 = $ws.Range("E142").Value2
$ws.Range("F141").Value2 = $ws.Range("F142").Value2
$ws.Range("G141").Value2 = $ws.Range("G142").Value2
$ws.Range("H141").Value2 = $ws.Range("H142").Value2
$ws.Range("I141").Value2 = $ws.Range("I142").Value2
$ws.Range("J141").Value2 = $ws.Range("J142").Value2
$ws.Range("K141").Value2 = $ws.Range("K142").Value2
$ws.Range("L141").Value2 = "Primera"
$ws.Range("M141").Value2 = 55
$ws.Range("N141").Value2 = 20000
$ws.Range("O141").Value2 = 20000
$ws.Range("P141").Value2 = 20000
$ws.Range("Q141").Value2 = "$/bandeja 15 kilos granel"
$ws.Range("R141").Value2 = "Región de O'Higgins"
$ws.Range("S141").Value2 = 1333
$ws.Range("T141").Value2 = 15
